$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Mejora" row (row 12) with its data.
$ws.Range("A12").Value = "Estimación y Esfuerzo"
$ws.Range("B12").Value = "Aclarar como trabajar con la columna de complejidad en la pestaña de Factor de complejidad. Especificando que los valores en esta columna no debían de ser cambiados."
$ws.Range("C12").Value = "Estimaciones y Planeación"
$ws.Range("D12").Value = "Actualización"
$ws.Range("E12").Value = "Aprobado "

# B12 needs to wrap its (longer) text, same as the other formatted cells.
$ws.Range("B12").Font.Name = "Arial"
$ws.Range("B12").Font.Size = 10
$ws.Range("B12").WrapText = $true

# Row grew taller to fit the new content.
$ws.Rows.Item(12).RowHeight = 61.5

# Move the active selection to C13, matching the saved view state.
$ws.Range("C13").Select()
